# Update the 국장_조선_분석 decision sheet with the latest model run figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (HD HYUNDAI MIPO)
$ws.Range("E2").Value = 59.1
$ws.Range("F2").Value = 1.59
$ws.Range("I2").Value = 60
$ws.Range("J2").Value = 73
$ws.Range("K2").Value = 58.8
$ws.Range("N2").Value = 66.04328690552585

# Row 3 (Hanwha Ocean)
$ws.Range("K3").Value = 53.8
$ws.Range("N3").Value = 66.04328690552585

# Row 4 (SamsungHvyInd)
$ws.Range("K4").Value = 53.2
$ws.Range("N4").Value = 66.04328690552585

# Row 5 (010140.KS)
$ws.Range("K5").Value = 53.2
$ws.Range("N5").Value = 66.04328690552585
